$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H header "Save", cloning the existing header style (bold,
# centered, bordered) from G1 so the new header matches the rest of row 1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# "Save" column values for rows 2-17
$values = @(0, 0, 0, 1, 0, 0, 0, 0, 1, 1, 0, 0, 0, 0, 0, 1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
